$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted at row 33, pushing all the
# subsequent records (old rows 33-82) down by one row (to rows 34-83).
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new record's data.
$ws.Cells.Item(33, 1).Value2  = 11
$ws.Cells.Item(33, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(33, 3).Value2  = "Bíobío"
$ws.Cells.Item(33, 4).Value2  = 44540
$ws.Cells.Item(33, 5).Value2  = 8
$ws.Cells.Item(33, 6).Value2  = 100112032
$ws.Cells.Item(33, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(33, 8).Value2  = "Huracán"
$ws.Cells.Item(33, 9).Value2  = "Primera"
$ws.Cells.Item(33, 10).Value2 = 190
$ws.Cells.Item(33, 11).Value2 = 6000
$ws.Cells.Item(33, 12).Value2 = 6500
$ws.Cells.Item(33, 13).Value2 = 6263
$ws.Cells.Item(33, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(33, 15).Value2 = "Región del Maule"
$ws.Cells.Item(33, 16).Value2 = 104
$ws.Cells.Item(33, 17).Value2 = 60
$ws.Cells.Item(33, 18).Value2 = "Hortaliza"
